# Add the student ID ("2220016") and course code ("E14-TP2") as two new
# centered title-page paragraphs, inserted between "Carl Trépanier" and
# "Exposé sur le langage Python" — matching paragraph's existing
# formatting (centered, Arial Rounded MT Bold, 00B0F0, 32pt).

$d = $word.ActiveDocument

$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Carl Trépanier", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0) | Out-Null
$anchor.InsertParagraphAfter() | Out-Null

$carlIndex = $anchor.Paragraphs(1).Index
$p1 = $d.Paragraphs($carlIndex + 1)
$p1.Range.Text = "2220016"

$p1.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs($carlIndex + 2)
$p2.Range.Text = "E14-TP2"

# Drop the stale lastRenderedPageBreak cache marker in front of
# "Utilisations de Python" (a rendering hint Word recomputes on layout;
# re-writing the run's text clears the stale cached marker).
$d.Content.Find.Execute("Utilisations de Python", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Utilisations de Python", 2) | Out-Null
